$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): update the Iron & steel value, clear the Non-metallic minerals value
$ws.Range("B3").Value = 1501515.362606632
$ws.Range("D3").Value = ""
$ws.Range("D3").Style = "Normal"

# Row 4 (Methanol): update the Chemicals value
$ws.Range("C4").Value = 170.2654775220202

# Row 5 (Ammonia): update the Chemicals value
$ws.Range("C5").Value = 4908.329972190175

# Row 7: rename "Other" -> "Biogas" and update its Non-metallic minerals value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 1414.767251341918

# Row 8 (new): re-add the "Other" row with its own value, matching the
# row-label formatting used by the other category cells in column A
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = ""
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = 1126.96699604008
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
